$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = "com.armedia.acm.services.mediaengine.model.MediaEngineBusinessProcessModel"
$ws.Range("E17").Value = "MediaEngineWorkFlow"
$ws.Range("E18").Value = "MediaEngineWorkFlow"
$ws.Range("C14").Value = "`$model: MediaEngineBusinessProcessModel"
